$wb = $excel.ActiveWorkbook

# 1. Remove the asset-category worksheets that contain no data rows
#    (only a header row): 船舶, 航空器, 現金, 債券, 基金受益憑證, 保險, 債權, 債務
$sheetsToRemove = @("船舶", "航空器", "現金", "債券", "基金受益憑證", "保險", "債權", "債務")
foreach ($name in $sheetsToRemove) {
    $wb.Worksheets.Item($name).Delete()
}

# 2. 股票 (stocks) sheet: normalise headers to the English dataframe schema
#    and append legislator_name / legislator_id / date columns.
$stockSheet = $wb.Worksheets.Item("股票")

$stockSheet.Cells.Item(1, 2).Value = "name"
$stockSheet.Cells.Item(1, 3).Value = "owner"
$stockSheet.Cells.Item(1, 4).Value = "quantity"
$stockSheet.Cells.Item(1, 5).Value = "face_value"
$stockSheet.Cells.Item(1, 6).Value = "currency"
$stockSheet.Cells.Item(1, 7).Value = "total"
$stockSheet.Cells.Item(1, 8).Value = "date"
$stockSheet.Cells.Item(1, 9).Value = "legislator_name"
$stockSheet.Cells.Item(1, 10).Value = "legislator_id"

$stockSheet.Cells.Item(2, 8).NumberFormat = "@"
$stockSheet.Cells.Item(2, 8).Value = "2011-11-22"
$stockSheet.Cells.Item(2, 9).Value = "葉宜津"
$stockSheet.Cells.Item(2, 10).Value = 855

$stockSheet.Cells.Item(3, 8).NumberFormat = "@"
$stockSheet.Cells.Item(3, 8).Value = "2011-11-22"
$stockSheet.Cells.Item(3, 9).Value = "葉宜津"
$stockSheet.Cells.Item(3, 10).Value = 855

$stockSheet.Cells.Item(4, 8).NumberFormat = "@"
$stockSheet.Cells.Item(4, 8).Value = "2011-11-22"
$stockSheet.Cells.Item(4, 9).Value = "葉宜津"
$stockSheet.Cells.Item(4, 10).Value = 855

# 3. 其他有價證券 (other securities) sheet: the trailing H column is entirely
#    empty - drop it so the sheet's dimension shrinks back to A:G.
$otherSecSheet = $wb.Worksheets.Item("其他有價證券")
$otherSecSheet.Columns.Item(8).Delete()
